$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$CENTER = -4108

function Set-DateCell($row, $text) {
    $c = $ws.Cells.Item($row, 1)
    $c.HorizontalAlignment = $CENTER
    $c.Value = "'" + $text
}

function Set-HoursCell($row, $hours) {
    $c = $ws.Cells.Item($row, 2)
    $c.HorizontalAlignment = $CENTER
    $c.Value = $hours
}

function Set-DescCell($row, $text) {
    $c = $ws.Cells.Item($row, 3)
    $c.WrapText = $true
    $c.Value = $text
}

# Seed new shared strings in target order (90..96)
Set-DateCell 53 'April 7, 2025'
Set-DescCell 53 '100% done with PlayAnalytics (now called PlayAnalytics 2.5 ready for user feedback and presentation) for PHASE 4: PlayAnalytics (Sports Management Portal) Presentation to Client which I uploaded to Repo (Implementation > PlayAnalytics 2.5 ready for user feedback and presentation) . Basically 3 of 3 features (finished the Profiles page, Prediction page, and Dashboard page) is finally done ready for user feedback and defense presentation'
Set-DescCell 55 'Made some minor tweaks with PlayAnalytics (now called PlayAnalytics 2.5.1 ready for user feedback and presentation) for PHASE 4: PlayAnalytics (Sports Management Portal) Presentation to Client which I uploaded to Repo (Implementation > PlayAnalytics 2.5.1 ready for user feedback and presentation) and created a zipped file of PlayAnalytics which is outside the Implementation folder so that it is ready for user feedback and presentation'
Set-DateCell 57 'April 10, 2025'
Set-DescCell 54 'Created a PlayAnalytics User Feedback Survey using Microsoft Forms for PHASE 4: PlayAnalytics (Sports Management Portal) Presentation to Client and shared link with 21 persons with deadline to gather responses by April 10, Thursday'
Set-DescCell 57 'Uploaded to Repo "PlayAnalytics User Feedback Survey(1-17).xlsx" under ReportsAndDocuments folder collecting user feedback from 17 out of 21 respondents for PHASE 4: PlayAnalytics (Sports Management Portal) Presentation to Client'
Set-DescCell 56 'Updated  file "README.md" and uploaded to Repo'
Set-DescCell 58 'Updated  file "README.md" and uploaded to Repo'

# Number of hours column
Set-HoursCell 53 2
Set-HoursCell 54 1
Set-HoursCell 55 1
Set-HoursCell 56 1
Set-HoursCell 57 1
Set-HoursCell 58 1

# Row heights to match wrapped-text rendering
$ws.Rows.Item(53).RowHeight = 102
$ws.Rows.Item(54).RowHeight = 68
$ws.Rows.Item(55).RowHeight = 102
$ws.Rows.Item(56).RowHeight = 17
$ws.Rows.Item(57).RowHeight = 68
$ws.Rows.Item(58).RowHeight = 17

# Update view/selection to the new bottom of the table
$ws.Range("A53").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 53
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C63").Select() | Out-Null

